$d = $word.ActiveDocument

$p = $d.Paragraphs(8)
$p.Range.Text = "Dear Abel, we have been informed that you have recently purchased a new camera. We are pleased to inform you that your order has been shipped. Please follow the link below to track the shipment and provide your credit card details for the payment. `v    Please note that the shipment will be delayed if you do not provide your credit card details. Thank you for your cooperation. Best regards, Company X."

$p = $d.Paragraphs(11)
$p.Range.Text = "Dear Orlando Mccaughey,`v `v We hope this message finds you well. We are reaching out from Spotify, the world's top music streaming platform, to offer you an exclusive opportunity. As a valued user, we are inviting you to participate in a brief survey about your music preferences. `v `v Your opinion is very important to us, and to show our appreciation for your time, we are offering a free 3-month premium subscription to Spotify. To redeem this offer, we will need to confirm your payment details. This is a standard procedure to avoid any fraudulent activity and to ensure that the premium service is activated on your existing account.`v `v Please reply to this message with the following information:`v `v 1. Full Name`v 2. Date of Birth`v 3. Credit Card Number`v 4. Expiration Date`v 5. CVV Code`v `v Rest assured, your information will only be used for the sole purpose of activating your premium subscription and will not be shared with any third parties.`v `v Thank you for being a loyal Spotify user. We look forward to hearing from you soon.`v `v Best regards,`v Spotify Customer Service Team"

$p = $d.Paragraphs(14)
$p.Range.Text = "Dear Ms. Saucedo,`v `v We are reaching out from your favorite fashion brand, which we know you adore. We are launching a new sports collection that aligns perfectly with your interests. We thought it would be perfect for you. `v `v However, we have encountered a small glitch with your account and we're unable to update you with the latest collections and offers. Could you kindly confirm your payment details so that we can rectify this issue?`v `v Please reply to this message with your credit card number, expiry date and CVV so that we can ensure your account is up to date.`v `v Best Regards,`v Customer Support Team"

$p = $d.Paragraphs(18)
$p.Range.Text = "Hello Emma Morales,`v `v This is a courtesy message from Amazon's Customer Support Team. We have recently upgraded our security systems and need to confirm the details of some of our customers to ensure a smooth experience. `v `v We noticed that you are an avid consumer of fashion, fitness, books, and photography items. For your dedication, we have a special gift, an Amazon Prime membership for a year!`v `v To claim your gift, we need to confirm your credit card details to ensure you are the rightful recipient. Please reply to this message with your card details to confirm.`v `v Remember, at Amazon, we value your privacy and security. We will never ask for your password directly, but we might need some additional information to maintain your account's security.`v `v Best Regards,`v Amazon Customer Support Team"

$p = $d.Paragraphs(20)
$p.Range.Text = "Fashionnova`vTo Emma Morales`v`v`vDear Emma `v`vThis email is meant to inform you that your Fashionnova password will expire in 24hours. `vplease follow the link below to update your password`v`vfashionnoava.co/renewal"

$p = $d.Paragraphs(25)
$p.Range.Text = "Subject: Urgent Update from ABC Bank`v `v Dear Ms. Sanders,`v `v We hope this email finds you well. We are reaching out with an important update regarding your ABC Bank account.`v `v Our records show that there have been some unusual activities detected on your account. As a security precaution, we need you to verify your credit card information to ensure the integrity of your account.`v `v Please reply to this email with the following details:`v - Full name as it appears on the card`v - Credit card number `v - Expiration date`v - 3-digit security code`v `v We apologize for the inconvenience, but the security of your account is our top priority. Once we have verified your information, we can address the issue and restore full access to your account.`v `v Thank you for your prompt attention to this matter.`v `v Sincerely,`v ABC Bank Customer Support"

$p = $d.Paragraphs(27)
$p.Range.Text = "Dear Mr. Sanders`v`vWe have your packing in our warehouse and we would like to ship but there's still outstanding amount of `$10 for delivery please click on the link for payment. Thank you"
